# se mejora parte de AP seleccion de metodo de pago No cuotas
# se mofifica data para smoke en preprod despues de trabajos de Oracle
# Se modifica testRun de Smoke

$wb = $excel.ActiveWorkbook

# --- DatosCuenta: update smoke QA data to new preprod data ---
$wsCuenta = $wb.Worksheets.Item("DatosCuenta")
$wsCuenta.Range("A2").Value = "SmokPreSiete"
$wsCuenta.Range("B2").Value = "SmokePreProdNamSiete"
$wsCuenta.Range("E2").Select()

# --- DatosMotor: move selection (testRun position) ---
$wsMotor = $wb.Worksheets.Item("DatosMotor")
$wsMotor.Activate()
$wsMotor.Range("E8").Select()

# --- DatosAP: move selection (testRun position) ---
$wsAP = $wb.Worksheets.Item("DatosAP")
$wsAP.Activate()
$wsAP.Range("H8").Select()

# Re-activate DatosAP as the final active sheet (matches tabSelected in DatosAP)
$wsAP.Activate()
